$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move the "MCT-1A-Tecnologia dos Materiais." class out of the Friday 10:40/11:30
# slots (F7/F8) and into the Tuesday 8:40/9:50 slots (C4/C6).
$ws.Range("C4").Value = "[-, 'MCT-1A-Tecnologia dos Materiais.']"
$ws.Range("C6").Value = "[-, 'MCT-1A-Tecnologia dos Materiais.']"
$ws.Range("F7").Value = "-"
$ws.Range("F8").Value = "-"
